# Update gh-pages to output generated at 456a3b4
#
# This script applies the scraped-data edit to the workbook:
#  - bumps a handful of "想去人数" (interested-count) numbers on sheet
#    "展览" (Exhibitions) and sheet "全部类型" (All types)
#  - inserts one brand-new event row ("杭州·排球少年only·春日校庆") into
#    both of those listings, just before the "杭州·创造力动漫游戏嘉年华1.0"
#    row, pushing every row after it down by one
#  - bumps the minimum ticket price for the 夏川里美 concert row on sheet
#    "演出" (Performances) and on its mirrored row in "全部类型"

$wb = $excel.ActiveWorkbook

function Set-FValues($ws, $map) {
    foreach ($row in $map.Keys) {
        $ws.Cells.Item($row, 6).Value = $map[$row]
    }
}

function Insert-NewEventRow($ws, $rowIndex, $lastRow) {
    # Push rows rowIndex..end down by one, matching Excel's native
    # "Insert Row" behaviour (formats/row heights shift with the data).
    $ws.Rows.Item($rowIndex).Insert()

    # The freshly inserted row inherits a blank/default style; pull the
    # real per-column formatting back from the row that is now directly
    # below our insertion point (it used to be the same physical row,
    # so its formatting is exactly what the new row should have).
    $below = $rowIndex + 1
    $ws.Range($ws.Cells.Item($below, 1), $ws.Cells.Item($below, 9)).Copy()
    $ws.Range($ws.Cells.Item($rowIndex, 1), $ws.Cells.Item($rowIndex, 9)).PasteSpecial(-4122)

    # Force text columns that look numeric/date-like to stay plain text
    # (matches every other row in the sheet - inlineStr, General style)
    # instead of being auto-parsed into a date serial by COM.
    $ws.Cells.Item($rowIndex, 2).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 2).Value = "2024-03-30"
    $ws.Cells.Item($rowIndex, 3).Value = "杭州·排球少年only·春日校庆"
    $ws.Cells.Item($rowIndex, 4).Value = "之江路149号 云栖培训基地"
    $ws.Cells.Item($rowIndex, 5).Value = "2024.03.30 10:00-03.31 17:00"
    $ws.Cells.Item($rowIndex, 6).Value = 1
    $ws.Cells.Item($rowIndex, 7).Value = 89
    $ws.Cells.Item($rowIndex, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81511"
    $ws.Cells.Item($rowIndex, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/qJrJxGUy1706581833764.jpeg"

    # Clean up the quotePrefix/number-format style artifact the forced
    # text-entry above left on column B, again by pulling the plain
    # format back from a neighbouring already-correct cell.
    $ws.Cells.Item($below, 2).Copy()
    $ws.Cells.Item($rowIndex, 2).PasteSpecial(-4122)

    # Column A is a manually-maintained running index (row number - 1),
    # not a formula/autofill, so it does NOT follow the row shift on its
    # own - every row from the insertion point down to the old last row
    # needs to be renumbered by +1 (the brand-new last row keeps the
    # index the old last row used to have).
    for ($r = $rowIndex; $r -le ($lastRow + 1); $r++) {
        $ws.Cells.Item($r, 1).Value = ($r - 1)
    }
}

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

Set-FValues $wsExpo @{
    3  = 7892
    5  = 939
    6  = 291
    9  = 92
    13 = 3153
    15 = 96
    19 = 458
    21 = 252
    22 = 224
    23 = 319
    26 = 108
    27 = 276
    31 = 502
    32 = 525
    33 = 23
}

Insert-NewEventRow $wsExpo 35 37

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$wsShows = $wb.Worksheets.Item("演出")
$wsShows.Cells.Item(6, 7).Value = 280

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) - same edits as "展览" plus the mirrored
# price bump from "演出", all offset by +2 rows (it also carries the two
# rows that live on "演出"/"本地生活" ahead of the exhibition listing).
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

Set-FValues $wsAll @{
    5  = 7893
    7  = 939
    8  = 291
    11 = 92
    16 = 3153
    18 = 96
    24 = 458
    26 = 252
    27 = 224
    28 = 319
    31 = 108
    32 = 276
    36 = 502
    37 = 525
    38 = 23
}

Insert-NewEventRow $wsAll 40 43

# After the insert, the 夏川里美 concert row (previously row 43) is now
# row 44; bump its min ticket price the same way as on "演出".
$wsAll.Cells.Item(44, 7).Value = 280
